$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("H33").Value = 5365.2856
$ws.Range("I33").Value = 156.3
$ws.Range("J33").Value = 10100.728
$ws.Range("K33").Value = 156.3
$ws.Range("L33").Value = 10100.728
$ws.Range("M33").Value = 72.69999999999999
$ws.Range("N33").Value = -10558.728
$ws.Range("H87").Value = 11602.774
$ws.Range("J87").Value = 11602.774
$ws.Range("L87").Value = 11602.774
$ws.Range("N87").Value = -14098.774
$ws.Range("H90").Value = 11602.774
$ws.Range("J90").Value = 11602.774
$ws.Range("L90").Value = 34808.322
$ws.Range("N90").Value = -47288.322
$ws.Range("H100").Value = 14494305
$ws.Range("J100").Value = 4753
$ws.Range("L100").Value = 4753
$ws.Range("N100").Value = -5835
$ws.Range("H137").Value = 1125.7273
$ws.Range("I137").Value = 1182.5625
$ws.Range("J137").Value = 1072.2354
$ws.Range("K137").Value = 3547.6875
$ws.Range("L137").Value = 3216.7062
$ws.Range("M137").Value = -997.6875
$ws.Range("N137").Value = -8316.706200000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1620.7894
$ws.Range("I2").Value = 621.7
$ws.Range("J2").Value = 2730.889
$ws.Range("K2").Value = 621.7
$ws.Range("L2").Value = 2730.889
$ws.Range("M2").Value = -508.7
$ws.Range("N2").Value = -2956.889
$ws.Range("H45").Value = 4422.25
$ws.Range("I45").Value = 5120.8
$ws.Range("J45").Value = 3258
$ws.Range("K45").Value = 5120.8
$ws.Range("L45").Value = 3258
$ws.Range("M45").Value = -4743.8
$ws.Range("N45").Value = -4012
$ws.Range("H110").Value = 632.5263
$ws.Range("I110").Value = 663.75
$ws.Range("J110").Value = 579
$ws.Range("K110").Value = 663.75
$ws.Range("L110").Value = 579
$ws.Range("M110").Value = 1381.25
$ws.Range("N110").Value = -4669
$ws.Range("H116").Value = 1620.7894
$ws.Range("I116").Value = 621.7
$ws.Range("J116").Value = 2730.889
$ws.Range("K116").Value = 621.7
$ws.Range("L116").Value = 2730.889
$ws.Range("M116").Value = 1672.3
$ws.Range("N116").Value = -7318.889

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1620.7894
$ws.Range("I3").Value = 621.7
$ws.Range("J3").Value = 2730.889
$ws.Range("K3").Value = 621.7
$ws.Range("L3").Value = 2730.889
$ws.Range("M3").Value = -507.7
$ws.Range("N3").Value = -2958.889
$ws.Range("H33").Value = 3900
$ws.Range("I33").Value = 1800
$ws.Range("K33").Value = 1800
$ws.Range("M33").Value = -1464
$ws.Range("H80").Value = 581.4375
$ws.Range("I80").Value = 1900.6666
$ws.Range("J80").Value = 277
$ws.Range("K80").Value = 1900.6666
$ws.Range("L80").Value = 277
$ws.Range("M80").Value = -902.6666
$ws.Range("N80").Value = -2273
$ws.Range("H83").Value = 581.4375
$ws.Range("I83").Value = 1900.6666
$ws.Range("J83").Value = 277
$ws.Range("K83").Value = 9503.333000000001
$ws.Range("L83").Value = 1385
$ws.Range("M83").Value = -4511.333000000001
$ws.Range("N83").Value = -11369
$ws.Range("H99").Value = 3303.5334
$ws.Range("I99").Value = 1361.4286
$ws.Range("J99").Value = 5002.875
$ws.Range("K99").Value = 1361.4286
$ws.Range("L99").Value = 5002.875
$ws.Range("M99").Value = 136.5714
$ws.Range("N99").Value = -7998.875
$ws.Range("H107").Value = 2331.7
$ws.Range("I107").Value = 1757.7142
$ws.Range("J107").Value = 3671
$ws.Range("K107").Value = 1757.7142
$ws.Range("L107").Value = 3671
$ws.Range("M107").Value = 162.2858000000001
$ws.Range("N107").Value = -7511

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 55557444
$ws.Range("I58").Value = 142858260
$ws.Range("J58").Value = 2381.818
$ws.Range("K58").Value = 142858260
$ws.Range("L58").Value = 2381.818
$ws.Range("M58").Value = -142858057
$ws.Range("N58").Value = -2787.818
$ws.Range("H88").Value = 8000
$ws.Range("I88").Value = 8000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 8000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -7594
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 8000
$ws.Range("I91").Value = 8000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 8000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -6596
$ws.Range("N91").ClearContents()
$ws.Range("H107").Value = 32258422
$ws.Range("I107").Value = 41667016
$ws.Range("J107").Value = 395.14285
$ws.Range("K107").Value = 41667016
$ws.Range("L107").Value = 395.14285
$ws.Range("M107").Value = -41665096
$ws.Range("N107").Value = -4235.14285
$ws.Range("H122").Value = 1473.1666
$ws.Range("I122").Value = 966.4
$ws.Range("J122").Value = 4007
$ws.Range("K122").Value = 2899.2
$ws.Range("L122").Value = 12021
$ws.Range("M122").Value = -449.1999999999998
$ws.Range("N122").Value = -16921
$ws.Range("H136").Value = 55557444
$ws.Range("I136").Value = 142858260
$ws.Range("J136").Value = 2381.818
$ws.Range("K136").Value = 428574780
$ws.Range("L136").Value = 7145.454000000001
$ws.Range("M136").Value = -428572230
$ws.Range("N136").Value = -12245.454

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H32").Value = 8390108
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 8390108
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 25170324
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -25170890
$ws.Range("H60").Value = 9009.111000000001
$ws.Range("I60").Value = 413.125
$ws.Range("J60").Value = 77777
$ws.Range("K60").Value = 1239.375
$ws.Range("L60").Value = 233331
$ws.Range("M60").Value = -988.375
$ws.Range("N60").Value = -233833
$ws.Range("H68").Value = 920.3871
$ws.Range("I68").Value = 767.1111
$ws.Range("J68").Value = 983.0909
$ws.Range("K68").Value = 2301.3333
$ws.Range("L68").Value = 2949.2727
$ws.Range("M68").Value = -1490.3333
$ws.Range("N68").Value = -4571.2727
$ws.Range("H71").Value = 920.3871
$ws.Range("I71").Value = 767.1111
$ws.Range("J71").Value = 983.0909
$ws.Range("K71").Value = 6903.9999
$ws.Range("L71").Value = 8847.8181
$ws.Range("M71").Value = -2847.9999
$ws.Range("N71").Value = -16959.8181
$ws.Range("H97").Value = 435.5
$ws.Range("I97").Value = 445.33334
$ws.Range("K97").Value = 1336.00002
$ws.Range("M97").Value = -840.0000199999999
$ws.Range("H98").Value = 1476.3125
$ws.Range("I98").Value = 1296.5
$ws.Range("J98").Value = 1502
$ws.Range("K98").Value = 3889.5
$ws.Range("L98").Value = 4506
$ws.Range("M98").Value = -2391.5
$ws.Range("N98").Value = -7502
$ws.Range("H125").Value = 8000
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H134").Value = 2344.5
$ws.Range("J134").Value = 6333.3335
$ws.Range("L134").Value = 19000.0005
$ws.Range("N134").Value = -29140.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 25000
$ws.Range("J46").Value = 25000
$ws.Range("L46").Value = 25000
$ws.Range("N46").Value = -25312
$ws.Range("H122").Value = 2380
$ws.Range("I122").Value = 2380
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7140
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4690
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 20834882
$ws.Range("J126").Value = 166666670
$ws.Range("L126").Value = 500000010
$ws.Range("N126").Value = -500004950

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 625.73914
$ws.Range("I22").Value = 560
$ws.Range("J22").Value = 749
$ws.Range("K22").Value = 560
$ws.Range("L22").Value = 749
$ws.Range("M22").Value = -265
$ws.Range("N22").Value = -1339
$ws.Range("H27").Value = 625.73914
$ws.Range("I27").Value = 560
$ws.Range("J27").Value = 749
$ws.Range("K27").Value = 560
$ws.Range("L27").Value = 749
$ws.Range("M27").Value = -453
$ws.Range("N27").Value = -963
$ws.Range("H55").Value = 564.73914
$ws.Range("I55").Value = 389.2857
$ws.Range("J55").Value = 837.6667
$ws.Range("K55").Value = 389.2857
$ws.Range("L55").Value = 837.6667
$ws.Range("M55").Value = -216.2857
$ws.Range("N55").Value = -1183.6667
$ws.Range("H61").Value = 1969.2
$ws.Range("I61").Value = 1183.4736
$ws.Range("J61").Value = 3326.3635
$ws.Range("K61").Value = 1183.4736
$ws.Range("L61").Value = 3326.3635
$ws.Range("M61").Value = -981.4736
$ws.Range("N61").Value = -3730.3635
$ws.Range("H113").Value = 1969.2
$ws.Range("I113").Value = 1183.4736
$ws.Range("J113").Value = 3326.3635
$ws.Range("K113").Value = 1183.4736
$ws.Range("L113").Value = 3326.3635
$ws.Range("M113").Value = 986.5264
$ws.Range("N113").Value = -7666.363499999999
$ws.Range("H136").Value = 75818.86
$ws.Range("I136").Value = 143640.58
$ws.Range("J136").Value = 7997.143
$ws.Range("K136").Value = 430921.74
$ws.Range("L136").Value = 23991.429
$ws.Range("M136").Value = -428371.74
$ws.Range("N136").Value = -29091.429

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H81").Value = 2914.697
$ws.Range("I81").Value = 2485.6875
$ws.Range("J81").Value = 3318.4707
$ws.Range("K81").Value = 4971.375
$ws.Range("L81").Value = 6636.9414
$ws.Range("M81").Value = -3910.375
$ws.Range("N81").Value = -8758.9414
$ws.Range("H84").Value = 2914.697
$ws.Range("I84").Value = 2485.6875
$ws.Range("J84").Value = 3318.4707
$ws.Range("K84").Value = 24856.875
$ws.Range("L84").Value = 33184.70699999999
$ws.Range("M84").Value = -19552.875
$ws.Range("N84").Value = -43792.70699999999
$ws.Range("H107").Value = 792.2222
$ws.Range("I107").Value = 755.2857
$ws.Range("K107").Value = 2265.8571
$ws.Range("M107").Value = -345.8571000000002
$ws.Range("H113").Value = 350.6875
$ws.Range("I113").Value = 353.33334
$ws.Range("J113").Value = 342.75
$ws.Range("K113").Value = 1060.00002
$ws.Range("L113").Value = 1028.25
$ws.Range("M113").Value = 1109.99998
$ws.Range("N113").Value = -5368.25
$ws.Range("H122").Value = 2225.25
$ws.Range("I122").Value = 2093
$ws.Range("J122").Value = 2357.5
$ws.Range("K122").Value = 6279
$ws.Range("L122").Value = 7072.5
$ws.Range("M122").Value = -3829
$ws.Range("N122").Value = -11972.5
